$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3 values
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 1

# Copy formatting of A2 (style index 1) down to A4:A5 to match the existing A column style
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Add new row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1

# Add new row 5
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 1
